$d = $word.ActiveDocument

function Add-Bullet($anchorPara, $text, $level) {
    $anchorPara.Range.InsertParagraphAfter() | Out-Null
    $np = $anchorPara.Next()
    $np.Range.Text = $text
    $np.Range.ListFormat.ListLevelNumber = $level
    return $np
}

# ---------------------------------------------------------------------------
# 1) Normalise paragraph 1's text (it was split over 3 runs: "I" + "mplement
#    ..." + " (do all this ...)"). Find/Replace matches across run
#    boundaries and rewrites it as a single run with the very same text.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Implement all common array operations that includes traversal like finding maximum, minimum, sum, product, linear search, reverse array. (do all this solution in single class and try to reuse the loop for function).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Implement all common array operations that includes traversal like finding maximum, minimum, sum, product, linear search, reverse array. (do all this solution in single class and try to reuse the loop for function).",
    2) | Out-Null

# NOTE: Paragraph handles in this host re-resolve lazily by their current
# document position, so we always walk strictly forward with `.Next()`
# immediately before using a handle, and never "look ahead" and cache a
# handle before an earlier insertion would shift its position.

$p1 = $d.Paragraphs.Item(1)          # Implement all common array operations...
$p2 = $p1.Next()                     # Write a code that swaps every alternate...

# ---------------------------------------------------------------------------
# 2) New sub-bullet right after "Write a code ...".
# ---------------------------------------------------------------------------
$p3 = Add-Bullet $p2 "if array = [10,20,30,40,50,60] => output array => [20,10,40,30,60,50]" 2

$p4 = $p3.Next()                     # If array = [10,20,30,40,50] => output array => [20,10,40,30,50]
$p5 = $p4.Next()                     # If array = [10,20,30,40,50,60] => output array => [20,10,40,30,60,50] (to be retargeted)

# ---------------------------------------------------------------------------
# 3) The old "If array = [10,20,30,40,50,60] => output array => ..." bullet
#    (originally split over 5 runs) gets reused/retargeted to a brand new
#    sentence. Find/Replace across its runs lands a single clean run.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "If array = [10,20,30,40,50,60] => output array => [20,10,40,30,60,50]",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Array = [ 10,20,25,20,10] then answer =>25.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 4) New sub-bullet right after the retargeted "...answer =>25." paragraph.
# ---------------------------------------------------------------------------
$cursor = Add-Bullet $p5 "Array = [ 20,30,40,45,40,30,20] then answer =>45." 2

# ---------------------------------------------------------------------------
# 5) "Find unique value in odd size array. ..." stays right where it always
#    was, immediately following the paragraph we just inserted. All brand
#    new material gets appended directly below it.
# ---------------------------------------------------------------------------
$cursor = $cursor.Next()             # Find unique value in odd size array. ...
$cursor = Add-Bullet $cursor "Find whether the values in an array are in unique number of occurrence or not" 1
$cursor = Add-Bullet $cursor "If array = [1,2,3,4] => returns false because 1,2,3,4 elements are repeating one time." 2
$cursor = Add-Bullet $cursor "If array = [1,2,2,3,3,3] => return true because every value has unique number of occurrence" 2
$cursor = Add-Bullet $cursor "Find duplicate value in array which contains 1 to n-1 values. You can also consider the same problem which doesn't contain 1 to n-1 values." 1
$cursor = Add-Bullet $cursor " If array = [1,3,2,4,3] => output => 3" 2
$cursor = Add-Bullet $cursor " If array = [1,2,3,2] => output => 2" 2
$cursor = Add-Bullet $cursor "Find Intersection of 2 sorted arrays. Intersection means the same elements both arrays can have" 1
$cursor = Add-Bullet $cursor "If array-1 = [1,2,5,7] and array-2 = [5,6,7] => output => [5,7]" 2
$cursor = Add-Bullet $cursor " If array-1 = [1,2,3,4,5] and array-2 = [5,6] => output => [5]" 2
$cursor = Add-Bullet $cursor "Find Intersection of 3 sorted arrays" 1
$cursor = Add-Bullet $cursor "Find pairs of indexes that gives sum equals to target." 1
$cursor = Add-Bullet $cursor "Index should be returned in sorted order" 2
$cursor = Add-Bullet $cursor "Example: - array => [1,2,3,4,5] and target => 6 Then answer => [ [ 0,4], [ 1,3]]" 2
$cursor = Add-Bullet $cursor "Find triplet of indexes that gives sum equals to target." 1
$cursor = Add-Bullet $cursor "sort an array with values 0,1." 1
$cursor = Add-Bullet $cursor "sort an array with values 0,1,2" 1

# ---------------------------------------------------------------------------
# 6) Apply the new paragraph spacing (1.5 line spacing, 12pt before) to
#    every paragraph in the document, including the trailing bookmark-only
#    paragraph.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $p.Range.ParagraphFormat.SpaceBefore = 12
    $p.Range.ParagraphFormat.LineSpacingRule = 1
}

Write-Output "Paragraphs: $($d.Paragraphs.Count)"
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    Write-Output "$i [$($p.Range.ListFormat.ListLevelNumber)]: $($p.Range.Text)"
}
